# Adds two new registrant rows (42 and 43) to the registration sheet,
# matching the rows submitted on 2025-02-24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 42 ----
$ws.Range("A42").Value = 41
$ws.Range("B42").Value = "Nguyễn Thị Phương Anh"
$ws.Range("C42").Value = "Bác sĩ"
$ws.Range("D42").Value = "Khoa Liên chuyên khoa"
$ws.Range("E42").Value = "Có"

# Phone number: force text so the leading zero survives.
$ws.Range("F42").NumberFormat = "@"
$ws.Range("F42").Value = "0363558688"
$ws.Range("F42").Style = "Normal"

$ws.Range("G42").Value = "Phuonganhhom@gmail.com"
$ws.Range("H42").Value = "2025-02-24 08:25:11"

# Date of birth: force text so it is not parsed into a date serial number.
$ws.Range("I42").NumberFormat = "@"
$ws.Range("I42").Value = "12/07/1993"
$ws.Range("I42").Style = "Normal"

$ws.Range("J42").Value = "Điều dưỡng"
$ws.Range("K42").Value = "Cao đẳng"

# ---- Row 43 ----
$ws.Range("A43").Value = 42
$ws.Range("B43").Value = "Phạm Thị Thu"
$ws.Range("C43").Value = "Bác sĩ"
$ws.Range("D43").Value = "Khoa Liên chuyên khoa"
$ws.Range("E43").Value = "Có"

# Phone number: force text so the leading zero survives.
$ws.Range("F43").NumberFormat = "@"
$ws.Range("F43").Value = "0983564922"
$ws.Range("F43").Style = "Normal"

$ws.Range("G43").Value = "bsthubvdonganh@gmail.com"
$ws.Range("H43").Value = "2025-02-24 08:26:08"
$ws.Range("I43").Value = "19/11/1983"
$ws.Range("J43").Value = "Bác sỹ"
$ws.Range("K43").Value = "ThS"
